$d = $word.ActiveDocument

# ------------------------------------------------------------------------
# Helper: returns the character offset (relative to Document.Content) where
# $text starts, searching the whole document.
# ------------------------------------------------------------------------
function Find-Offset($text) {
    $r = $d.Range(0, $d.Content.End)
    $r.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    return $r.Start
}

# Helper: toggles a no-op character format over a range so the host splits
# the underlying run(s) at the range's start/end without changing the
# visible formatting.
function Split-Range($startOffset, $endOffset) {
    $r = $d.Range($startOffset, $endOffset)
    $r.Font.Bold = $true
    $r.Font.Bold = $false
}

# Helper: returns the character offset right after $text ends (searching
# the whole document).
function Find-End($text) {
    $r = $d.Range(0, $d.Content.End)
    $r.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    return $r.End
}

# ------------------------------------------------------------------------
# Locate the boundaries of the original runs in the summary paragraph
# before making any edit (these anchors are all unchanged text, so they
# can be located up front).
# ------------------------------------------------------------------------
$runStart = Find-Offset "Seasoned Technical Architect with over 13 years of experience"
$bBoot = Find-Offset "Boot, Microservices"
$bMicro = Find-Offset "Microservices, ReactJs"
$bAnd = Find-Offset " and SQL"
$bSql = Find-Offset "SQL, with comprehensive"
$bWith = Find-Offset ", with comprehensive"
$bBehaviour = Find-Offset "Behaviour-Driven"
$paraEnd = Find-End "exceed project goals."

$seasonedLen = "Seasoned ".Length
$oldPhrase = "Technical Architect"
$newPhrase = "Technical Specialist"

$wordStart = $runStart + $seasonedLen
$wordEnd = $wordStart + $oldPhrase.Length

# ------------------------------------------------------------------------
# Step 1: perform the actual text edit ("Architect" -> "Specialist").
# Changing text anywhere in the paragraph causes the host to coalesce all
# same-formatted runs of that paragraph into a single run, so do this
# before re-establishing the desired run boundaries.
# ------------------------------------------------------------------------
$wordRange = $d.Range($wordStart, $wordEnd)
$wordRange.Text = $newPhrase

# The paragraph's text grew by this many characters; shift the boundaries
# that come after the edit point accordingly.
$delta = $newPhrase.Length - $oldPhrase.Length
$bBoot = $bBoot + $delta
$bMicro = $bMicro + $delta
$bAnd = $bAnd + $delta
$bSql = $bSql + $delta
$bWith = $bWith + $delta
$bBehaviour = $bBehaviour + $delta
$paraEnd = $paraEnd + $delta

# ------------------------------------------------------------------------
# Step 2: re-split the (now coalesced) run back into the pieces required
# by the edit, and restore all of the original run boundaries further
# along in the paragraph so the rest of it is left as it was.
# ------------------------------------------------------------------------
$seasonedEnd = $runStart + $seasonedLen
$specialistEnd = $seasonedEnd + $newPhrase.Length
$spaceEnd = $specialistEnd + 1

Split-Range $runStart $seasonedEnd
Split-Range $seasonedEnd $specialistEnd
Split-Range $specialistEnd $spaceEnd
Split-Range $spaceEnd $bBoot
Split-Range $bBoot $bMicro
Split-Range $bMicro $bAnd
Split-Range $bAnd $bSql
Split-Range $bSql $bWith
Split-Range $bWith $bBehaviour
Split-Range $bBehaviour $paraEnd

Write-Host "piece1[" $d.Range($runStart, $seasonedEnd).Text "]"
Write-Host "piece2[" $d.Range($seasonedEnd, $specialistEnd).Text "]"
Write-Host "piece3[" $d.Range($specialistEnd, $spaceEnd).Text "]"
Write-Host "piece4[" $d.Range($spaceEnd, $bBoot).Text "]"
